# "haciendo union entre empresa y contacto"
# Fill in missing contact emails (with mailto hyperlinks) and append two more
# rappi contacts (rows 14-15) that round out the rappi company group.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Add-EmailHyperlink($cellRef, $email, $formatDonor) {
    $range = $ws.Range($cellRef)
    $ws.Hyperlinks.Add($range, "mailto:$email", [Type]::Missing, [Type]::Missing, $email)
    # Hyperlinks.Add applies Excel's built-in blue/underlined "Hyperlink" style;
    # this workbook keeps its existing per-row style on email cells instead, so
    # restore that by pasting the (unstyled) format from a sibling cell.
    $ws.Range($formatDonor).Copy()
    $range.PasteSpecial(-4122) | Out-Null
}

# G5: "ggdfgfdgfd" -> "juanluis@gmail.com"
Add-EmailHyperlink "G5" "juanluis@gmail.com" "H5"

# G10: "correo1" -> "correo1@gmail.com"
Add-EmailHyperlink "G10" "correo1@gmail.com" "H10"

# G11: "correo3" -> "correo3@gmail.com"
Add-EmailHyperlink "G11" "correo3@gmail.com" "H11"

# G12: "correo4" -> "gggg@gmail.com"
Add-EmailHyperlink "G12" "gggg@gmail.com" "H12"

# Row 14 (new): rappi6 contact
$ws.Range("A14").Value = "rappi6"
$ws.Range("B14").Value = "rappi.com"
$ws.Range("C14").Value = "delivery"
$ws.Range("D14").Value = "mateus2"
$ws.Range("E14").Value = "uribe"
$ws.Range("F14").Value = 543545
$ws.Range("H14").Value = "precualification"
Add-EmailHyperlink "G14" "correo6@gmail.com" "H14"

# Row 15 (new): rappi7 contact
$ws.Range("A15").Value = "rappi7"
$ws.Range("B15").Value = "rappi.com"
$ws.Range("C15").Value = "delivery"
$ws.Range("D15").Value = "mateus2"
$ws.Range("E15").Value = "uribe"
$ws.Range("F15").Value = 543545
$ws.Range("H15").Value = "precualification"
Add-EmailHyperlink "G15" "correo7@gmail.com" "H15"

# Move active selection to F14, matching the workbook's last-edited cell
$ws.Range("F14").Select()
